$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift the old "blank template" rows 35/36 down to 37/38 ---
# (done manually, cell-by-cell, rather than via Rows.Insert, so that no
# stray formatted-but-empty cells get introduced in columns that merely
# carry a column-level style, e.g. column E)
$ws.Range("A37").Value = "LeetCode"
$ws.Range("A37").HorizontalAlignment = -4131
$ws.Range("B37").Value = "Anna"
$ws.Range("B37").HorizontalAlignment = -4131
$ws.Range("C37").Value = "Easy"

$ws.Range("A38").Value = "LeetCode"
$ws.Range("A38").HorizontalAlignment = -4131
$ws.Range("B38").Value = "Stephan"
$ws.Range("B38").HorizontalAlignment = -4131
$ws.Range("C38").Value = "Easy"

# --- New row 35: LeetCode / Anna / Easy / 706. Design HashMap ---
$ws.Range("A35").Value = "LeetCode"
$ws.Range("A35").HorizontalAlignment = -4131
$ws.Range("B35").Value = "Anna"
$ws.Range("B35").HorizontalAlignment = -4131
$ws.Range("C35").Value = "Easy"
$ws.Range("D35").WrapText = $true
$ws.Range("D35").Value = "706. Design HashMap"
$ws.Rows("35").RowHeight = 14

# --- New row 36: LeetCode / Stephan / Easy / 706. Design HashMap / 2020/12/23 / Completed ---
$ws.Range("A36").Value = "LeetCode"
$ws.Range("A36").HorizontalAlignment = -4131
$ws.Range("B36").Value = "Stephan"
$ws.Range("B36").HorizontalAlignment = -4131
$ws.Range("C36").Value = "Easy"
$ws.Range("D36").WrapText = $true
$ws.Range("D36").Value = "706. Design HashMap"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "2020/12/23"
$ws.Range("G36").Value = "Completed"
$ws.Rows("36").RowHeight = 14

# --- Update the saved selection to match the workbook view ---
$ws.Range("D40").Select()
